# Add daily reset for detailed messages:
# append 8 new incident rows (213-220) to the detailed messages log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value that must stay plain text even though it
# looks like a date/time/duration, without leaving a non-default
# cell style behind (matches the unstyled inlineStr cells used by
# every other row in the sheet).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value2 = $value
    $range.Style = "Normal"
}

$rows = @(
    @{ Row = 213; A = "WC49 P5H"; B = "La cámara no detecta Busbar";          C = "2024-06-11"; D = "11:27:21"; E = "Mañana"; F = "11:27:24"; G = "0:00:03"; H = "-0.00 minutos" },
    @{ Row = 214; A = "WC49 P5H"; B = "No coloca bien el sealling";           C = "2024-06-11"; D = "11:27:37"; E = "Mañana"; F = "11:27:41"; G = "0:00:04"; H = "0.07 minutos" },
    @{ Row = 215; A = "WC49 P5H"; B = "La cámara no detecta Busbar";          C = "2024-06-11"; D = "11:37:49"; E = "Mañana"; F = "11:38:06"; G = "0:00:17"; H = "-0.00 minutos" },
    @{ Row = 216; A = "WC49 P5H"; B = "Screw K30 no lo detecta puesto";       C = "2024-06-11"; D = "11:53:14"; E = "Mañana"; F = $null;      G = $null;     H = "-0.01 minutos" },
    @{ Row = 217; A = "WC49 P5H"; B = "Etiquetadora";                        C = "2024-06-11"; D = "11:55:22"; E = "Mañana"; F = "11:55:25"; G = "0:00:03"; H = "-0.00 minutos" },
    @{ Row = 218; A = "WC49 P5H"; B = "Power atascado en prensa, cuesta sacar"; C = "2024-06-11"; D = "11:55:34"; E = "Mañana"; F = "11:55:35"; G = "0:00:01"; H = "0.05 minutos" },
    @{ Row = 219; A = "WC49 P5H"; B = "La cámara no detecta Top Cover";       C = "2024-06-11"; D = "12:10:03"; E = "Mañana"; F = "12:10:05"; G = "0:00:02"; H = "-0.00 minutos" },
    @{ Row = 220; A = "WC49 P5H"; B = "La cámara no detecta Top Cover";       C = "2024-06-11"; D = "12:12:11"; E = "Mañana"; F = "12:12:12"; G = "0:00:01"; H = "-0.01 minutos" }
)

foreach ($r in $rows) {
    $n = $r.Row

    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = $r.B

    Set-TextValue $ws.Range("C$n") $r.C
    Set-TextValue $ws.Range("D$n") $r.D

    $ws.Range("E$n").Value = $r.E

    if ($r.F) {
        Set-TextValue $ws.Range("F$n") $r.F
    }
    if ($r.G) {
        Set-TextValue $ws.Range("G$n") $r.G
    }

    $ws.Range("H$n").Value = $r.H
}
